$d = $word.ActiveDocument

# The "In Attendance" list had "Rahat Nafees" followed by a trailing
# end-of-paragraph run containing just a space, and then a separate
# run containing the literal text "(Late)". Word had therefore split a
# single " (Late)" annotation across two adjacent runs that share the
# exact same run formatting (rStyle "eop", Aptos/Segoe UI fonts).
#
# Re-typing/replacing across that run boundary with identical
# replacement text collapses the two runs that made up " (Late)" into
# a single run, matching the cleaned-up markup.
$d.Content.Find.Execute(" (Late)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " (Late)", 2) | Out-Null

$d.Save()
